# Generate Report for Archive
# Reassign the three in-flight localization files (06edefed, 0a3cc689, d3e7fab8)
# to rows 7-9 on each sheet (Overview, zh-cn, de-de), all now "In Translation".
# Row 10 (d8e6a821) is unaffected.

$wb = $excel.ActiveWorkbook

# ---- Sheet "Overview" ----
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A7").Value = "06edefed-d24c-4db8-8c01-fbb61e353865.md"
$ws.Range("B7").Value = "In Translation"
$ws.Range("C7").Value = "In Translation"

$ws.Range("A8").Value = "0a3cc689-f94c-48a0-9726-46cc3c34de71.md"
$ws.Range("B8").Value = "In Translation"
$ws.Range("C8").Value = "In Translation"

$ws.Range("A9").Value = "d3e7fab8-0350-4ca7-86dd-e36c0d26afb4.md"
$ws.Range("B9").Value = "In Translation"
$ws.Range("C9").Value = "In Translation"

# ---- Sheet "zh-cn" ----
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A7").Value = "06edefed-d24c-4db8-8c01-fbb61e353865.md"
$ws.Range("B7").Value = "In Translation"
$ws.Range("C7").Value = "06edefed-d24c-4db8-8c01-fbb61e353865.59e821c7f4c743536980b9867460bd7f51270e3a.zh-cn.xlf"
$ws.Range("D7").Value = "2016-03-11 02:29:00"

$ws.Range("A8").Value = "0a3cc689-f94c-48a0-9726-46cc3c34de71.md"
$ws.Range("B8").Value = "In Translation"
$ws.Range("C8").Value = "0a3cc689-f94c-48a0-9726-46cc3c34de71.d3acb8505bf652ffae3bb1ebd63913790dd893bf.zh-cn.xlf"
$ws.Range("D8").Value = "2016-03-11 02:29:00"

$ws.Range("A9").Value = "d3e7fab8-0350-4ca7-86dd-e36c0d26afb4.md"
$ws.Range("B9").Value = "In Translation"
$ws.Range("C9").Value = "d3e7fab8-0350-4ca7-86dd-e36c0d26afb4.a9567d8361ef552a0252e4f39417c927a83e4a86.zh-cn.xlf"
$ws.Range("D9").Value = "2016-03-11 02:15:37"

# ---- Sheet "de-de" ----
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A7").Value = "06edefed-d24c-4db8-8c01-fbb61e353865.md"
$ws.Range("B7").Value = "In Translation"
$ws.Range("C7").Value = "06edefed-d24c-4db8-8c01-fbb61e353865.59e821c7f4c743536980b9867460bd7f51270e3a.de-de.xlf"
$ws.Range("D7").Value = "2016-03-11 02:29:07"

$ws.Range("A8").Value = "0a3cc689-f94c-48a0-9726-46cc3c34de71.md"
$ws.Range("B8").Value = "In Translation"
$ws.Range("C8").Value = "0a3cc689-f94c-48a0-9726-46cc3c34de71.d3acb8505bf652ffae3bb1ebd63913790dd893bf.de-de.xlf"
$ws.Range("D8").Value = "2016-03-11 02:29:07"

$ws.Range("A9").Value = "d3e7fab8-0350-4ca7-86dd-e36c0d26afb4.md"
$ws.Range("B9").Value = "In Translation"
$ws.Range("C9").Value = "d3e7fab8-0350-4ca7-86dd-e36c0d26afb4.a9567d8361ef552a0252e4f39417c927a83e4a86.de-de.xlf"
$ws.Range("D9").Value = "2016-03-11 02:17:03"
